# Applies the cell-value updates described in the diff for
# output/B_vector_diag.xlsx (column A binary vector) to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column A.
$updates = @{
    8   = 0
    34  = 0
    50  = 0
    55  = 0
    61  = 1
    63  = 0
    66  = 0
    68  = 0
    107 = 1
    118 = 1
    131 = 1
    132 = 0
    147 = 0
    156 = 1
    159 = 0
    167 = 0
    173 = 1
    175 = 1
    176 = 1
    199 = 0
    220 = 0
    222 = 0
    261 = 1
    295 = 1
    296 = 1
    300 = 0
    321 = 1
    328 = 0
    349 = 0
    390 = 1
    393 = 0
    400 = 1
    409 = 0
    418 = 0
    498 = 0
    547 = 0
    724 = 1
    761 = 0
    763 = 1
    797 = 0
    799 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
